# Auto-generated edit script: apply numeric corrections to Tonberry_Profits (multi-sheet) workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 11 (hunk 0, G=5533)
$ws.Range("H11").Value = 373
$ws.Range("I11").Value = 373
$ws.Range("K11").Value = 373
$ws.Range("M11").Value = -233

# row 18 (hunk 1, G=5471)
$ws.Range("H18").Value = 13691.52
$ws.Range("I18").Value = 11164.286
$ws.Range("J18").Value = 16908
$ws.Range("K18").Value = 11164.286
$ws.Range("L18").Value = 16908
$ws.Range("M18").Value = -10880.286
$ws.Range("N18").Value = -17476

# row 33 (hunk 2, G=5512)
$ws.Range("H33").Value = 195.83333
$ws.Range("I33").Value = 240.90909
$ws.Range("K33").Value = 240.90909
$ws.Range("M33").Value = -11.90908999999999

# row 81 (hunk 3, G=10637)
$ws.Range("H81").Value = 38000
$ws.Range("J81").Value = 38000
$ws.Range("L81").Value = 38000
$ws.Range("N81").Value = -39996

# row 84 (hunk 4, G=10637)
$ws.Range("H84").Value = 38000
$ws.Range("J84").Value = 38000
$ws.Range("L84").Value = 114000
$ws.Range("N84").Value = -123984

# row 116 (hunk 5, G=27778)
$ws.Range("H116").Value = 14868.7
$ws.Range("I116").Value = 51000
$ws.Range("J116").Value = 5835.875
$ws.Range("K116").Value = 51000
$ws.Range("L116").Value = 5835.875
$ws.Range("M116").Value = -47558
$ws.Range("N116").Value = -12719.875

# row 131 (hunk 6, G=36108)
$ws.Range("H131").Value = 1563.6666
$ws.Range("I131").Value = 668.1818
$ws.Range("J131").Value = 4026.25
$ws.Range("K131").Value = 2004.5454
$ws.Range("L131").Value = 12078.75
$ws.Range("M131").Value = 3035.4546
$ws.Range("N131").Value = -22158.75

# row 136 (hunk 7, G=42164)
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# row 137 (hunk 8, G=44013)
$ws.Range("H137").Value = 1785.875
$ws.Range("I137").Value = 1493.5
$ws.Range("J137").Value = 1883.3334
$ws.Range("K137").Value = 4480.5
$ws.Range("L137").Value = 5650.0002
$ws.Range("M137").Value = -1930.5
$ws.Range("N137").Value = -10750.0002

# row 138 (hunk 9, G=44169)
$ws.Range("H138").Value = 2216.75
$ws.Range("I138").Value = 2236.5386
$ws.Range("K138").Value = 6709.6158
$ws.Range("M138").Value = -1569.6158

# row 140 (hunk 10, G=42459)
$ws.Range("H140").Value = 47427.285
$ws.Range("J140").Value = 47427.285
$ws.Range("L140").Value = 47427.285
$ws.Range("N140").Value = -57787.285

$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk 11, G=44147)
$ws.Range("H32").Value = 3446.0784
$ws.Range("I32").Value = 2143.1052
$ws.Range("K32").Value = 2143.1052
$ws.Range("M32").Value = -1856.1052

# row 61 (hunk 12, G=43999)
$ws.Range("H61").Value = 6276.3335
$ws.Range("I61").Value = 4459.4
$ws.Range("J61").Value = 8547.5
$ws.Range("K61").Value = 4459.4
$ws.Range("L61").Value = 8547.5
$ws.Range("M61").Value = -4247.4
$ws.Range("N61").Value = -8971.5

# row 63 (hunk 13, G=12528)
$ws.Range("H63").Value = 8795
$ws.Range("J63").Value = 7992.3335
$ws.Range("L63").Value = 7992.3335
$ws.Range("N63").Value = -9364.333500000001

# row 66 (hunk 14, G=12528)
$ws.Range("H66").Value = 8795
$ws.Range("J66").Value = 7992.3335
$ws.Range("L66").Value = 39961.6675
$ws.Range("N66").Value = -46825.6675

# row 74 (hunk 15, G=44000)
$ws.Range("H74").Value = 4365.0625
$ws.Range("I74").Value = 4352.857
$ws.Range("K74").Value = 4352.857
$ws.Range("M74").Value = -3478.857

# row 77 (hunk 16, G=44000)
$ws.Range("H77").Value = 4365.0625
$ws.Range("I77").Value = 4352.857
$ws.Range("K77").Value = 21764.285
$ws.Range("M77").Value = -17396.285

# row 97 (hunk 17, G=19941)
$ws.Range("H97").Value = 1183.3334
$ws.Range("I97").Value = 1183.3334
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1183.3334
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -687.3334
$ws.Range("N97").ClearContents()

# row 102 (hunk 18, G=19945)
$ws.Range("H102").Value = 1036.8334
$ws.Range("I102").Value = 997.4
$ws.Range("K102").Value = 997.4
$ws.Range("M102").Value = 624.6

# row 132 (hunk 19, G=43997)
$ws.Range("H132").Value = 1961.0344
$ws.Range("I132").Value = 1195.3889
$ws.Range("J132").Value = 3213.9092
$ws.Range("K132").Value = 3586.1667
$ws.Range("L132").Value = 9641.7276
$ws.Range("M132").Value = -1056.1667
$ws.Range("N132").Value = -14701.7276

# row 136 (hunk 20, G=43999)
$ws.Range("H136").Value = 6276.3335
$ws.Range("I136").Value = 4459.4
$ws.Range("J136").Value = 8547.5
$ws.Range("K136").Value = 13378.2
$ws.Range("L136").Value = 25642.5
$ws.Range("M136").Value = -10828.2
$ws.Range("N136").Value = -30742.5

$ws = $wb.Worksheets.Item("BSM")
# row 26 (hunk 21, G=19535)
$ws.Range("H26").Value = 40000
$ws.Range("I26").Value = 40000
$ws.Range("K26").Value = 40000
$ws.Range("M26").Value = -39708

# row 96 (hunk 22, G=19525)
$ws.Range("H96").Value = 12750
$ws.Range("I96").Value = 12750
$ws.Range("K96").Value = 12750
$ws.Range("M96").Value = -10004

# row 134 (hunk 23, G=43998)
$ws.Range("H134").Value = 7414.2583
$ws.Range("I134").Value = 8416.583000000001
$ws.Range("J134").Value = 3977.7144
$ws.Range("K134").Value = 25249.749
$ws.Range("L134").Value = 11933.1432
$ws.Range("M134").Value = -22714.749
$ws.Range("N134").Value = -17003.1432

$ws = $wb.Worksheets.Item("CRP")
# row 5 (hunk 24, G=1893)
$ws.Range("H5").Value = 2452
$ws.Range("J5").Value = 2452
$ws.Range("L5").Value = 2452
$ws.Range("N5").Value = -2676

# row 31 (hunk 25, G=44023)
$ws.Range("H31").Value = 2097.5144
$ws.Range("I31").Value = 928.3461
$ws.Range("J31").Value = 5475.1113
$ws.Range("K31").Value = 928.3461
$ws.Range("L31").Value = 5475.1113
$ws.Range("M31").Value = -633.3461
$ws.Range("N31").Value = -6065.1113

# row 34 (hunk 26, G=44023)
$ws.Range("H34").Value = 2097.5144
$ws.Range("I34").Value = 928.3461
$ws.Range("J34").Value = 5475.1113
$ws.Range("K34").Value = 928.3461
$ws.Range("L34").Value = 5475.1113
$ws.Range("M34").Value = -726.3461
$ws.Range("N34").Value = -5879.1113

# row 58 (hunk 27, G=44021)
$ws.Range("H58").Value = 1368.4286
$ws.Range("I58").Value = 1251
$ws.Range("J58").Value = 1579.8
$ws.Range("K58").Value = 1251
$ws.Range("L58").Value = 1579.8
$ws.Range("M58").Value = -1048
$ws.Range("N58").Value = -1985.8

# row 69 (hunk 28, G=11911)
$ws.Range("H69").Value = 220201
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# row 72 (hunk 29, G=11911)
$ws.Range("H72").Value = 220201
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# row 122 (hunk 30, G=36196)
$ws.Range("H122").Value = 1328.5172
$ws.Range("I122").Value = 1327.5264
$ws.Range("J122").Value = 1330.4
$ws.Range("K122").Value = 3982.5792
$ws.Range("L122").Value = 3991.2
$ws.Range("M122").Value = -1532.5792
$ws.Range("N122").Value = -8891.200000000001

# row 132 (hunk 31, G=44019)
$ws.Range("H132").Value = 1919.6842
$ws.Range("I132").Value = 917.0417
$ws.Range("K132").Value = 2751.1251
$ws.Range("M132").Value = -221.1251000000002

# row 134 (hunk 32, G=44020)
$ws.Range("H134").Value = 1025.4615
$ws.Range("I134").Value = 1010.9167
$ws.Range("K134").Value = 3032.7501
$ws.Range("M134").Value = -497.7501000000002

# row 136 (hunk 33, G=44021)
$ws.Range("H136").Value = 1368.4286
$ws.Range("I136").Value = 1251
$ws.Range("J136").Value = 1579.8
$ws.Range("K136").Value = 3753
$ws.Range("L136").Value = 4739.4
$ws.Range("M136").Value = -1203
$ws.Range("N136").Value = -9839.4

$ws = $wb.Worksheets.Item("CUL")
# row 6 (hunk 34, G=4639)
$ws.Range("H6").Value = 63
$ws.Range("I6").Value = 63
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 189
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -76
$ws.Range("N6").ClearContents()

# row 131 (hunk 35, G=36060)
$ws.Range("H131").Value = 13908995
$ws.Range("J131").Value = 24828.207
$ws.Range("L131").Value = 74484.621
$ws.Range("N131").Value = -84564.621

$ws = $wb.Worksheets.Item("GSM")
# row 102 (hunk 36, G=36169)
$ws.Range("H102").Value = 2217.5454
$ws.Range("I102").Value = 2712.7778
$ws.Range("J102").Value = 1874.6923
$ws.Range("K102").Value = 2712.7778
$ws.Range("L102").Value = 1874.6923
$ws.Range("M102").Value = -1090.7778
$ws.Range("N102").Value = -5118.6923

# row 136 (hunk 37, G=42218)
$ws.Range("H136").Value = 9666.333000000001
$ws.Range("J136").Value = 9666.333000000001
$ws.Range("L136").Value = 28998.999
$ws.Range("N136").Value = -34098.999

$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 38, G=36249)
$ws.Range("H7").Value = 4205.6523
$ws.Range("I7").Value = 2273.8
$ws.Range("J7").Value = 5691.6924
$ws.Range("K7").Value = 2273.8
$ws.Range("L7").Value = 5691.6924
$ws.Range("M7").Value = -2161.8
$ws.Range("N7").Value = -5915.6924

# row 40 (hunk 39, G=36248)
$ws.Range("H40").Value = 5384.8076
$ws.Range("I40").Value = 3066.3333
$ws.Range("J40").Value = 7372.0713
$ws.Range("K40").Value = 3066.3333
$ws.Range("L40").Value = 7372.0713
$ws.Range("M40").Value = -2930.3333
$ws.Range("N40").Value = -7644.0713

# row 100 (hunk 40, G=19995)
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

# row 126 (hunk 41, G=36249)
$ws.Range("H126").Value = 4205.6523
$ws.Range("I126").Value = 2273.8
$ws.Range("J126").Value = 5691.6924
$ws.Range("K126").Value = 6821.400000000001
$ws.Range("L126").Value = 17075.0772
$ws.Range("M126").Value = -4351.400000000001
$ws.Range("N126").Value = -22015.0772

# row 136 (hunk 42, G=44060)
$ws.Range("H136").Value = 3378.8438
$ws.Range("I136").Value = 2655.125
$ws.Range("J136").Value = 5550
$ws.Range("K136").Value = 7965.375
$ws.Range("L136").Value = 16650
$ws.Range("M136").Value = -5415.375
$ws.Range("N136").Value = -21750

$ws = $wb.Worksheets.Item("WVR")
# row 81 (hunk 43, G=12596)
$ws.Range("H81").Value = 489
$ws.Range("I81").Value = 489
$ws.Range("K81").Value = 978
$ws.Range("M81").Value = 83

# row 84 (hunk 44, G=12596)
$ws.Range("H84").Value = 489
$ws.Range("I84").Value = 489
$ws.Range("K84").Value = 4890
$ws.Range("M84").Value = 414

# row 100 (hunk 45, G=19981)
$ws.Range("H100").Value = 350
$ws.Range("I100").Value = 350
$ws.Range("K100").Value = 700
$ws.Range("M100").Value = -159

# row 113 (hunk 46, G=27752)
$ws.Range("H113").Value = 796.9231
$ws.Range("J113").Value = 1520
$ws.Range("L113").Value = 4560
$ws.Range("N113").Value = -8900

# row 122 (hunk 47, G=36208)
$ws.Range("H122").Value = 35631.305
$ws.Range("I122").Value = 57387.785
$ws.Range("J122").Value = 1787.8889
$ws.Range("K122").Value = 172163.355
$ws.Range("L122").Value = 5363.6667
$ws.Range("M122").Value = -169713.355
$ws.Range("N122").Value = -10263.6667

# row 123 (hunk 48, G=34127)
$ws.Range("H123").Value = 47328.668
$ws.Range("J123").Value = 47328.668
$ws.Range("L123").Value = 47328.668
$ws.Range("N123").Value = -57128.668

# row 136 (hunk 49, G=44031)
$ws.Range("H136").Value = 3346.7
$ws.Range("I136").Value = 3048.5881
$ws.Range("J136").Value = 3736.5386
$ws.Range("K136").Value = 9145.764299999999
$ws.Range("L136").Value = 11209.6158
$ws.Range("M136").Value = -6595.764299999999
$ws.Range("N136").Value = -16309.6158

